$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 84
$ws1.Range("F3").Value = 11871
$ws1.Range("F4").Value = 14
$ws1.Range("F6").Value = 354
$ws1.Range("F8").Value = 11794
$ws1.Range("F13").Value = 1777
$ws1.Range("F14").Value = 5845
$ws1.Range("F16").Value = 3535
$ws1.Range("F17").Value = 187

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 84
$ws4.Range("F4").Value = 3
$ws4.Range("F5").Value = 11871
$ws4.Range("F6").Value = 14
$ws4.Range("F8").Value = 3
$ws4.Range("F9").Value = 354
$ws4.Range("F11").Value = 11794
$ws4.Range("F16").Value = 1777
$ws4.Range("F18").Value = 5845
$ws4.Range("F20").Value = 3535
$ws4.Range("F21").Value = 187
